$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold + border) from H1 into I1:J1, then set header labels
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data columns I (I0) and J (IF) for rows 2-87
$data = @{
    2 = @(8, 8)
    3 = @(6, 6)
    4 = @(5, 5)
    5 = @(7, 8)
    6 = @(10, 10)
    7 = @(8, 8)
    8 = @(8, 8)
    9 = @(8, 8)
    10 = @(7, 7)
    11 = @(7, 7)
    12 = @(7, 7)
    13 = @(5, 5)
    14 = @(5, 5)
    15 = @(6, 6)
    16 = @(7, 7)
    17 = @(6, 6)
    18 = @(7, 7)
    19 = @(6, 6)
    20 = @(6, 6)
    21 = @(7, 7)
    22 = @(8, 8)
    23 = @(6, 6)
    24 = @(8, 8)
    25 = @(6, 6)
    26 = @(7, 7)
    27 = @(7, 7)
    28 = @(6, 6)
    29 = @(4, 4)
    30 = @(8, 8)
    31 = @(9, 9)
    32 = @(6, 7)
    33 = @(7, 7)
    34 = @(8, 8)
    35 = @(8, 8)
    36 = @(7, 7)
    37 = @(8, 9)
    38 = @(7, 7)
    39 = @(7, 7)
    40 = @(8, 8)
    41 = @(6, 6)
    42 = @(6, 6)
    43 = @(7, 7)
    44 = @(9, 9)
    45 = @(6, 6)
    46 = @(3, 4)
    47 = @(7, 8)
    48 = @(8, 8)
    49 = @(8, 9)
    50 = @(7, 7)
    51 = @(5, 5)
    52 = @(12, 12)
    53 = @(8, 8)
    54 = @(6, 7)
    55 = @(9, 9)
    56 = @(6, 7)
    57 = @(6, 6)
    58 = @(5, 5)
    59 = @(8, 8)
    60 = @(6, 6)
    61 = @(8, 8)
    62 = @(4, 5)
    63 = @(9, 9)
    64 = @(7, 7)
    65 = @(8, 8)
    66 = @(6, 6)
    67 = @(8, 8)
    68 = @(7, 7)
    69 = @(7, 7)
    70 = @(8, 8)
    71 = @(6, 6)
    72 = @(6, 6)
    73 = @(9, 9)
    74 = @(7, 7)
    75 = @(7, 7)
    76 = @(7, 7)
    77 = @(8, 8)
    78 = @(7, 7)
    79 = @(6, 6)
    80 = @(7, 7)
    81 = @(6, 7)
    82 = @(9, 9)
    83 = @(7, 7)
    84 = @(8, 8)
    85 = @(8, 8)
    86 = @(5, 5)
    87 = @(6, 6)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item([int]$r, 9).Value = $vals[0]
    $ws.Cells.Item([int]$r, 10).Value = $vals[1]
}
